# Update row 4 values on the active worksheet (Cleveland/Columbus/Birmingham projections)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 8926
$ws.Range("B4").Value = 2859
$ws.Range("C4").Value = 18748
$ws.Range("D4").Value = 2.772265838509316
$ws.Range("E4").Value = 0.8878881987577639
$ws.Range("F4").Value = 5.822628571428573
$ws.Range("G4").Value = 0.4164993788819875
$ws.Range("H4").Value = 0.2248447204968944
$ws.Range("I4").Value = 0.6099037142857143
$ws.Range("J4").Value = 0.1273685714285715
$ws.Range("K4").Value = 0.09181428571428571
$ws.Range("L4").Value = 0.1663270857142857
$ws.Range("M4").Value = 8148
$ws.Range("N4").Value = 1506
$ws.Range("O4").Value = 26457
$ws.Range("P4").Value = 2.645766233766234
$ws.Range("Q4").Value = 0.489025974025974
$ws.Range("R4").Value = 8.589992837662338
$ws.Range("S4").Value = 0.3294981818181818
$ws.Range("T4").Value = 0.1228571428571428
$ws.Range("U4").Value = 0.6873832597402595
$ws.Range("V4").Value = 0.1622889142857143
$ws.Range("W4").Value = 0.09775619047619047
$ws.Range("X4").Value = 0.2961675474285713
